$wb = $excel.ActiveWorkbook

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Cells.Item(1,6).Value = "Recoil_E(MeV)"
$ws.Cells.Item(2,6).Value = 303.520972972973
$ws.Cells.Item(3,6).Value = 303.520972972973
$ws.Cells.Item(4,6).Value = 303.520972972973
$ws.Cells.Item(5,6).Value = 303.520972972973
$ws.Cells.Item(6,6).Value = 303.520972972973
$ws.Cells.Item(7,6).Value = 300.6575675675676
$ws.Cells.Item(8,6).Value = 300.6575675675676
$ws.Cells.Item(9,6).Value = 300.6575675675676
$ws.Cells.Item(10,6).Value = 300.6575675675676
$ws.Cells.Item(11,6).Value = 300.6575675675676
$ws.Cells.Item(12,6).Value = 300.6575675675676
$ws.Cells.Item(13,6).Value = 297.7941621621622
$ws.Cells.Item(14,6).Value = 297.7941621621622
$ws.Cells.Item(15,6).Value = 297.7941621621622
$ws.Cells.Item(16,6).Value = 297.7941621621622
$ws.Cells.Item(17,6).Value = 294.93075675675675
$ws.Cells.Item(18,6).Value = 294.93075675675675
$ws.Cells.Item(19,6).Value = 294.93075675675675
$ws.Cells.Item(20,6).Value = 294.93075675675675
$ws.Cells.Item(21,6).Value = 294.93075675675675
$ws.Cells.Item(22,6).Value = 292.06735135135136
$ws.Cells.Item(23,6).Value = 292.06735135135136
$ws.Cells.Item(24,6).Value = 292.06735135135136
$ws.Cells.Item(25,6).Value = 292.06735135135136
$ws.Cells.Item(26,6).Value = 292.06735135135136
$ws.Cells.Item(27,6).Value = 289.203945945946
$ws.Cells.Item(28,6).Value = 289.203945945946
$ws.Cells.Item(29,6).Value = 289.203945945946
$ws.Cells.Item(30,6).Value = 289.203945945946
$ws.Cells.Item(31,6).Value = 286.3405405405406
$ws.Cells.Item(32,6).Value = 286.3405405405406
$ws.Cells.Item(33,6).Value = 286.3405405405406
$ws.Cells.Item(34,6).Value = 283.47713513513514
$ws.Cells.Item(35,6).Value = 283.47713513513514
$ws.Cells.Item(36,6).Value = 283.47713513513514
$ws.Cells.Item(37,6).Value = 280.61372972972976
$ws.Cells.Item(38,6).Value = 280.61372972972976
$ws.Cells.Item(39,6).Value = 280.61372972972976
$ws.Cells.Item(40,6).Value = 277.7503243243243
$ws.Cells.Item(41,6).Value = 277.7503243243243
$ws.Cells.Item(42,6).Value = 274.886918918919
$ws.Cells.Item(43,6).Value = 274.886918918919
$ws.Name = "0um_420MeV"

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Cells.Item(1,6).Value = "Recoil_E(MeV)"
$ws.Cells.Item(2,6).Value = 245.1891891891892
$ws.Cells.Item(3,6).Value = 245.1891891891892
$ws.Cells.Item(4,6).Value = 242.91891891891893
$ws.Cells.Item(5,6).Value = 242.91891891891893
$ws.Cells.Item(6,6).Value = 242.91891891891893
$ws.Cells.Item(7,6).Value = 242.91891891891893
$ws.Cells.Item(8,6).Value = 242.91891891891893
$ws.Cells.Item(9,6).Value = 240.64864864864865
$ws.Cells.Item(10,6).Value = 240.64864864864865
$ws.Cells.Item(11,6).Value = 240.64864864864865
$ws.Cells.Item(12,6).Value = 240.64864864864865
$ws.Cells.Item(13,6).Value = 240.64864864864865
$ws.Cells.Item(14,6).Value = 240.64864864864865
$ws.Cells.Item(15,6).Value = 238.3783783783784
$ws.Cells.Item(16,6).Value = 238.3783783783784
$ws.Cells.Item(17,6).Value = 236.1081081081081
$ws.Cells.Item(18,6).Value = 236.1081081081081
$ws.Cells.Item(19,6).Value = 236.1081081081081
$ws.Cells.Item(20,6).Value = 233.83783783783784
$ws.Cells.Item(21,6).Value = 233.83783783783784
$ws.Cells.Item(22,6).Value = 229.2972972972973
$ws.Cells.Item(23,6).Value = 229.2972972972973
$ws.Cells.Item(24,6).Value = 229.2972972972973
$ws.Cells.Item(25,6).Value = 227.02702702702703
$ws.Cells.Item(26,6).Value = 227.02702702702703
$ws.Cells.Item(27,6).Value = 222.48648648648648
$ws.Name = "9um_333MeV"

# ---- Sheet 4 ----
$ws = $wb.Worksheets.Item(4)
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Cells.Item(1,6).Value = "Recoil_E(MeV)"
$ws.Cells.Item(2,6).Value = 264.0740720720721
$ws.Cells.Item(3,6).Value = 264.0740720720721
$ws.Cells.Item(4,6).Value = 264.0740720720721
$ws.Cells.Item(5,6).Value = 264.0740720720721
$ws.Cells.Item(6,6).Value = 264.0740720720721
$ws.Cells.Item(7,6).Value = 261.6060900900901
$ws.Cells.Item(8,6).Value = 261.6060900900901
$ws.Cells.Item(9,6).Value = 261.6060900900901
$ws.Cells.Item(10,6).Value = 261.6060900900901
$ws.Cells.Item(11,6).Value = 261.6060900900901
$ws.Cells.Item(12,6).Value = 259.13810810810816
$ws.Cells.Item(13,6).Value = 259.13810810810816
$ws.Cells.Item(14,6).Value = 259.13810810810816
$ws.Cells.Item(15,6).Value = 259.13810810810816
$ws.Cells.Item(16,6).Value = 259.13810810810816
$ws.Cells.Item(17,6).Value = 256.67012612612615
$ws.Cells.Item(18,6).Value = 256.67012612612615
$ws.Cells.Item(19,6).Value = 256.67012612612615
$ws.Cells.Item(20,6).Value = 254.20214414414417
$ws.Cells.Item(21,6).Value = 254.20214414414417
$ws.Cells.Item(22,6).Value = 254.20214414414417
$ws.Cells.Item(23,6).Value = 254.20214414414417
$ws.Cells.Item(24,6).Value = 251.7341621621622
$ws.Cells.Item(25,6).Value = 251.7341621621622
$ws.Cells.Item(26,6).Value = 251.7341621621622
$ws.Cells.Item(27,6).Value = 251.7341621621622
$ws.Cells.Item(28,6).Value = 251.7341621621622
$ws.Cells.Item(29,6).Value = 249.2661801801802
$ws.Cells.Item(30,6).Value = 249.2661801801802
$ws.Cells.Item(31,6).Value = 249.2661801801802
$ws.Cells.Item(32,6).Value = 246.79819819819824
$ws.Cells.Item(33,6).Value = 246.79819819819824
$ws.Cells.Item(34,6).Value = 246.79819819819824
$ws.Cells.Item(35,6).Value = 244.33021621621623
$ws.Cells.Item(36,6).Value = 244.33021621621623
$ws.Name = "6um_362MeV"

# ---- Sheet 5 ----
$ws = $wb.Worksheets.Item(5)
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Cells.Item(1,6).Value = "Recoil_E(MeV)"
$ws.Cells.Item(2,6).Value = 285.2292252252252
$ws.Cells.Item(3,6).Value = 285.2292252252252
$ws.Cells.Item(4,6).Value = 285.2292252252252
$ws.Cells.Item(5,6).Value = 285.2292252252252
$ws.Cells.Item(6,6).Value = 285.2292252252252
$ws.Cells.Item(7,6).Value = 282.5635315315315
$ws.Cells.Item(8,6).Value = 282.5635315315315
$ws.Cells.Item(9,6).Value = 282.5635315315315
$ws.Cells.Item(10,6).Value = 282.5635315315315
$ws.Cells.Item(11,6).Value = 282.5635315315315
$ws.Cells.Item(12,6).Value = 279.8978378378378
$ws.Cells.Item(13,6).Value = 279.8978378378378
$ws.Cells.Item(14,6).Value = 279.8978378378378
$ws.Cells.Item(15,6).Value = 279.8978378378378
$ws.Cells.Item(16,6).Value = 279.8978378378378
$ws.Cells.Item(17,6).Value = 279.8978378378378
$ws.Cells.Item(18,6).Value = 277.2321441441441
$ws.Cells.Item(19,6).Value = 277.2321441441441
$ws.Cells.Item(20,6).Value = 277.2321441441441
$ws.Cells.Item(21,6).Value = 274.5664504504504
$ws.Cells.Item(22,6).Value = 274.5664504504504
$ws.Cells.Item(23,6).Value = 274.5664504504504
$ws.Cells.Item(24,6).Value = 274.5664504504504
$ws.Cells.Item(25,6).Value = 271.9007567567568
$ws.Cells.Item(26,6).Value = 271.9007567567568
$ws.Cells.Item(27,6).Value = 271.9007567567568
$ws.Cells.Item(28,6).Value = 269.2350630630631
$ws.Cells.Item(29,6).Value = 269.2350630630631
$ws.Cells.Item(30,6).Value = 266.5693693693694
$ws.Cells.Item(31,6).Value = 266.5693693693694
$ws.Cells.Item(32,6).Value = 266.5693693693694
$ws.Cells.Item(33,6).Value = 263.9036756756757
$ws.Cells.Item(34,6).Value = 263.9036756756757
$ws.Cells.Item(35,6).Value = 263.9036756756757
$ws.Cells.Item(36,6).Value = 258.5722882882883
$ws.Cells.Item(37,6).Value = 258.5722882882883
$ws.Cells.Item(38,6).Value = 250.5752072072072
$ws.Cells.Item(39,6).Value = 250.5752072072072
$ws.Name = "3um_391MeV"

# ---- Sheet 6 ----
$ws = $wb.Worksheets.Item(6)
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Cells.Item(1,6).Value = "Recoil_E(MeV)"
$ws.Cells.Item(2,6).Value = 223.0997837837838
$ws.Cells.Item(3,6).Value = 223.0997837837838
$ws.Cells.Item(4,6).Value = 223.0997837837838
$ws.Cells.Item(5,6).Value = 221.03404504504505
$ws.Cells.Item(6,6).Value = 221.03404504504505
$ws.Cells.Item(7,6).Value = 221.03404504504505
$ws.Cells.Item(8,6).Value = 221.03404504504505
$ws.Cells.Item(9,6).Value = 221.03404504504505
$ws.Cells.Item(10,6).Value = 218.96830630630632
$ws.Cells.Item(11,6).Value = 218.96830630630632
$ws.Cells.Item(12,6).Value = 218.96830630630632
$ws.Cells.Item(13,6).Value = 218.96830630630632
$ws.Cells.Item(14,6).Value = 218.96830630630632
$ws.Cells.Item(15,6).Value = 216.90256756756756
$ws.Cells.Item(16,6).Value = 216.90256756756756
$ws.Cells.Item(17,6).Value = 216.90256756756756
$ws.Cells.Item(18,6).Value = 214.83682882882883
$ws.Cells.Item(19,6).Value = 214.83682882882883
$ws.Cells.Item(20,6).Value = 214.83682882882883
$ws.Cells.Item(21,6).Value = 212.77109009009007
$ws.Cells.Item(22,6).Value = 212.77109009009007
$ws.Cells.Item(23,6).Value = 210.70535135135137
$ws.Cells.Item(24,6).Value = 210.70535135135137
$ws.Cells.Item(25,6).Value = 208.6396126126126
$ws.Cells.Item(26,6).Value = 208.6396126126126
$ws.Name = "12um_303MeV"
